# Edit script generated to match target diff
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): update F column (想去人数) values ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2,6).Value2 = 230
$ws1.Cells.Item(6,6).Value2 = 88
$ws1.Cells.Item(8,6).Value2 = 383
$ws1.Cells.Item(9,6).Value2 = 4624
$ws1.Cells.Item(10,6).Value2 = 4624
$ws1.Cells.Item(13,6).Value2 = 1083
$ws1.Cells.Item(14,6).Value2 = 597
$ws1.Cells.Item(15,6).Value2 = 4119
$ws1.Cells.Item(16,6).Value2 = 156
$ws1.Cells.Item(17,6).Value2 = 159
$ws1.Cells.Item(19,6).Value2 = 199
$ws1.Cells.Item(20,6).Value2 = 3402
$ws1.Cells.Item(24,6).Value2 = 2954
$ws1.Cells.Item(25,6).Value2 = 121
$ws1.Cells.Item(26,6).Value2 = 119
$ws1.Cells.Item(29,6).Value2 = 176
$ws1.Cells.Item(32,6).Value2 = 46
$ws1.Cells.Item(36,6).Value2 = 5338
$ws1.Cells.Item(37,6).Value2 = 736
$ws1.Cells.Item(41,6).Value2 = 25
$ws1.Cells.Item(42,6).Value2 = 1051
$ws1.Cells.Item(43,6).Value2 = 437
$ws1.Cells.Item(45,6).Value2 = 1937
$ws1.Cells.Item(46,6).Value2 = 296
$ws1.Cells.Item(48,6).Value2 = 687
$ws1.Cells.Item(49,6).Value2 = 823

# --- Sheet "演出" (sheet2): insert new row 5 with new event, shift rest down ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(5).Insert()
# copy column-A number style (bold/centered/bordered index style) from row above into new row 5
$ws2.Cells.Item(4,1).Copy()
$ws2.Cells.Item(5,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Cells.Item(5,1).Value2 = 4
$ws2.Cells.Item(5,2).Formula = "'2024-06-21"
$ws2.Cells.Item(5,3).Formula = "'北京·【梦境重启，经典共鸣】 —— 神级日漫金曲演唱会"
$ws2.Cells.Item(5,4).Formula = "'复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$ws2.Cells.Item(5,5).Formula = "'2024.06.21 19:30-06.21 21:00"
$ws2.Cells.Item(5,6).Value2 = 0
$ws2.Cells.Item(5,7).Value2 = 116.5
$ws2.Cells.Item(5,8).Formula = "'https://show.bilibili.com/platform/detail.html?id=87128"
$ws2.Cells.Item(5,9).Formula = "'//i0.hdslb.com/bfs/openplatform/202406/SdvF1jOT1717765398224.jpeg"

# --- Sheet "演出": row that was row21 (Marcin Patrzalek) is now row22 after the insert; bump F (want-to-go count) 716 -> 717 ---
$ws2.Cells.Item(22,6).Value2 = 717

# --- Sheet "全部类型" (sheet4): update F column (想去人数) values ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(5,6).Value2 = 230
$ws4.Cells.Item(8,6).Value2 = 88
$ws4.Cells.Item(10,6).Value2 = 383
$ws4.Cells.Item(11,6).Value2 = 4624
$ws4.Cells.Item(12,6).Value2 = 4624
$ws4.Cells.Item(13,6).Value2 = 37
$ws4.Cells.Item(18,6).Value2 = 1083
$ws4.Cells.Item(19,6).Value2 = 597
$ws4.Cells.Item(20,6).Value2 = 4119
$ws4.Cells.Item(21,6).Value2 = 156
$ws4.Cells.Item(22,6).Value2 = 159
$ws4.Cells.Item(23,6).Value2 = 199
$ws4.Cells.Item(24,6).Value2 = 3402
$ws4.Cells.Item(25,6).Value2 = 2954
$ws4.Cells.Item(26,6).Value2 = 121
$ws4.Cells.Item(27,6).Value2 = 119
$ws4.Cells.Item(29,6).Value2 = 176
$ws4.Cells.Item(37,6).Value2 = 5338
$ws4.Cells.Item(39,6).Value2 = 736
$ws4.Cells.Item(44,6).Value2 = 1051
$ws4.Cells.Item(45,6).Value2 = 437
$ws4.Cells.Item(47,6).Value2 = 1937
$ws4.Cells.Item(49,6).Value2 = 687
$ws4.Cells.Item(50,6).Value2 = 823
